$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-35 down to 10-36.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with data (same constant columns as the
# other rows in this block, with new values for D, J, K, L, M, P).
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44175
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112029
$ws.Range("G9").Value = "Orégano"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("N9").Value = '$/docena de atados'
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 4000
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = "Hortaliza"

# Give the new row's date cell the same number format used by the other
# date cells in column D (style index 2 in the original workbook).
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
